$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row from diff hunk starting at original line 3470
$ws.Range("H57").Value = 21489.75
$ws.Range("J57").Value = 8319.666999999999
$ws.Range("L57").Value = 24959.001
$ws.Range("N57").Value = -25957.001
# row from diff hunk starting at original line 4419
$ws.Range("I76").Value = 5999.6665
$ws.Range("K76").Value = 5999.6665
$ws.Range("M76").Value = -5684.6665
# row from diff hunk starting at original line 4569
$ws.Range("I79").Value = 5999.6665
$ws.Range("K79").Value = 5999.6665
$ws.Range("M79").Value = -4907.6665
# row from diff hunk starting at original line 5426
$ws.Range("H96").Value = 402.15384
$ws.Range("I96").Value = 293.45456
$ws.Range("K96").Value = 880.36368
$ws.Range("M96").Value = 492.63632
# row from diff hunk starting at original line 6886
$ws.Range("H125").Value = 7409656
$ws.Range("I125").Value = 1841
$ws.Range("J125").Value = 11113564
$ws.Range("K125").Value = 16569
$ws.Range("L125").Value = 100022076
$ws.Range("M125").Value = -14109
$ws.Range("N125").Value = -100026996
# row from diff hunk starting at original line 6987
$ws.Range("H127").Value = 2626.2856
$ws.Range("I127").Value = 477.1
$ws.Range("K127").Value = 1431.3
$ws.Range("M127").Value = 3528.7
# row from diff hunk starting at original line 7088
$ws.Range("H129").Value = 2246.2104
$ws.Range("I129").Value = 1187.375
$ws.Range("J129").Value = 3016.2727
$ws.Range("K129").Value = 3562.125
$ws.Range("L129").Value = 9048.8181
$ws.Range("M129").Value = 1437.875
$ws.Range("N129").Value = -19048.8181
# row from diff hunk starting at original line 7293
$ws.Range("H133").Value = 52271.91
$ws.Range("J133").Value = 52271.91
$ws.Range("L133").Value = 52271.91
$ws.Range("N133").Value = -62391.91
# row from diff hunk starting at original line 7544
$ws.Range("H138").Value = 5585.206
$ws.Range("J138").Value = 11256.5
$ws.Range("L138").Value = 33769.5
$ws.Range("N138").Value = -44049.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row from diff hunk starting at original line 9145
$ws.Range("H29").Value = 5010
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row from diff hunk starting at original line 15672
$ws.Range("H22").Value = 279
$ws.Range("I22").Value = 227.5
$ws.Range("K22").Value = 227.5
$ws.Range("M22").Value = -54.5
# row from diff hunk starting at original line 19676
$ws.Range("H105").Value = 1374.875
$ws.Range("I105").Value = 1099.3636
$ws.Range("K105").Value = 1099.3636
$ws.Range("M105").Value = 647.6364000000001
# row from diff hunk starting at original line 20497
$ws.Range("H122").Value = 75000
$ws.Range("J122").Value = 75000
$ws.Range("L122").Value = 75000
$ws.Range("N122").Value = -84800
# row from diff hunk starting at original line 21229
$ws.Range("H137").Value = 58332.75
$ws.Range("J137").Value = 58332.75
$ws.Range("L137").Value = 58332.75
$ws.Range("N137").Value = -68532.75
# row from diff hunk starting at original line 21327
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row from diff hunk starting at original line 21764
$ws.Range("H6").Value = 169758
$ws.Range("I6").Value = 333533
$ws.Range("K6").Value = 333533
$ws.Range("M6").Value = -333420
# row from diff hunk starting at original line 24027
$ws.Range("H52").Value = 65700.664
$ws.Range("I52").Value = 65741
$ws.Range("J52").Value = 65695.625
$ws.Range("K52").Value = 65741
$ws.Range("L52").Value = 65695.625
$ws.Range("M52").Value = -65447
$ws.Range("N52").Value = -66283.625
# row from diff hunk starting at original line 24309
$ws.Range("H58").Value = 6434.2104
$ws.Range("I58").Value = 5339.636
$ws.Range("J58").Value = 7939.25
$ws.Range("K58").Value = 5339.636
$ws.Range("L58").Value = 7939.25
$ws.Range("M58").Value = -5136.636
$ws.Range("N58").Value = -8345.25
# row from diff hunk starting at original line 25687
$ws.Range("H86").Value = 100008
$ws.Range("J86").Value = 100008
$ws.Range("L86").Value = 100008
$ws.Range("N86").Value = -102254
# row from diff hunk starting at original line 25831
$ws.Range("H89").Value = 100008
$ws.Range("J89").Value = 100008
$ws.Range("L89").Value = 500040
$ws.Range("N89").Value = -511272
# row from diff hunk starting at original line 27908
$ws.Range("H132").Value = 5325.905
$ws.Range("I132").Value = 5279
$ws.Range("J132").Value = 5377.5
$ws.Range("K132").Value = 15837
$ws.Range("L132").Value = 16132.5
$ws.Range("M132").Value = -13307
$ws.Range("N132").Value = -21192.5
# row from diff hunk starting at original line 28009
$ws.Range("H134").Value = 872533.0600000001
$ws.Range("I134").Value = 558290.75
$ws.Range("K134").Value = 1674872.25
$ws.Range("M134").Value = -1672337.25
# row from diff hunk starting at original line 28110
$ws.Range("H136").Value = 6434.2104
$ws.Range("I136").Value = 5339.636
$ws.Range("J136").Value = 7939.25
$ws.Range("K136").Value = 16018.908
$ws.Range("L136").Value = 23817.75
$ws.Range("M136").Value = -13468.908
$ws.Range("N136").Value = -28917.75
# row from diff hunk starting at original line 28260
$ws.Range("H139").Value = 73745
$ws.Range("J139").Value = 73745
$ws.Range("L139").Value = 73745
$ws.Range("N139").Value = -84025

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row from diff hunk starting at original line 31228
$ws.Range("H56").Value = 7062.25
$ws.Range("I56").Value = 7062.25
$ws.Range("K56").Value = 7062.25
$ws.Range("M56").Value = -6532.25
# row from diff hunk starting at original line 32780
$ws.Range("H87").Value = 12998.75
$ws.Range("J87").Value = 14995
$ws.Range("L87").Value = 44985
$ws.Range("N87").Value = -47481
# row from diff hunk starting at original line 32930
$ws.Range("H90").Value = 12998.75
$ws.Range("J90").Value = 14995
$ws.Range("L90").Value = 134955
$ws.Range("N90").Value = -147435
# row from diff hunk starting at original line 34145
$ws.Range("H114").Value = 90911304
$ws.Range("I114").Value = 333333540
$ws.Range("J114").Value = 2972.125
$ws.Range("K114").Value = 1000000620
$ws.Range("L114").Value = 8916.375
$ws.Range("M114").Value = -999997366
$ws.Range("N114").Value = -15424.375

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row from diff hunk starting at original line 39460
$ws.Range("H80").Value = 7685.5713
$ws.Range("I80").Value = 6059.8
$ws.Range("K80").Value = 6059.8
$ws.Range("M80").Value = -5061.8
# row from diff hunk starting at original line 39607
$ws.Range("H83").Value = 7685.5713
$ws.Range("I83").Value = 6059.8
$ws.Range("K83").Value = 30299
$ws.Range("M83").Value = -25307
# row from diff hunk starting at original line 40269
$ws.Range("H97").Value = 994.73334
$ws.Range("I97").Value = 910.6316
$ws.Range("J97").Value = 1140
$ws.Range("K97").Value = 910.6316
$ws.Range("L97").Value = 1140
$ws.Range("M97").Value = -414.6316
$ws.Range("N97").Value = -2132
# row from diff hunk starting at original line 41056
$ws.Range("H113").Value = 634716.5600000001
$ws.Range("I113").Value = 1669667
$ws.Range("K113").Value = 1669667
$ws.Range("M113").Value = -1667497

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row from diff hunk starting at original line 44694
$ws.Range("H46").Value = 5820.8945
$ws.Range("I46").Value = 5578.357
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 5578.357
$ws.Range("L46").Value = 6500
$ws.Range("M46").Value = -5390.357
$ws.Range("N46").Value = -6876
# row from diff hunk starting at original line 48394
$ws.Range("H122").Value = 557878.1
$ws.Range("I122").Value = 2316
$ws.Range("K122").Value = 6948
$ws.Range("M122").Value = -4498
# row from diff hunk starting at original line 49172
$ws.Range("H138").Value = 66126.336
$ws.Range("J138").Value = 66126.336
$ws.Range("L138").Value = 66126.336
$ws.Range("N138").Value = -76406.336

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row from diff hunk starting at original line 54243
$ws.Range("H100").Value = 673.65216
$ws.Range("I100").Value = 675.05
$ws.Range("K100").Value = 1350.1
$ws.Range("M100").Value = -809.0999999999999
# row from diff hunk starting at original line 55511
$ws.Range("H126").Value = 2473.923
$ws.Range("I126").Value = 1901
$ws.Range("K126").Value = 5703
$ws.Range("M126").Value = -3233
